$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
# Row 33
$ws.Range("H33").Value = 197.5
$ws.Range("I33").Value = 217.625
$ws.Range("K33").Value = 217.625
$ws.Range("M33").Value = 11.375
# Row 45
$ws.Range("H45").Value = 3728
$ws.Range("J45").Value = 4657.4287
$ws.Range("L45").Value = 13972.2861
$ws.Range("N45").Value = -14356.2861
# Row 54
$ws.Range("H54").Value = 504000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 96
$ws.Range("H96").Value = 1510
$ws.Range("I96").Value = 1158.2
$ws.Range("K96").Value = 3474.6
$ws.Range("M96").Value = -2101.6
# Row 100
$ws.Range("H100").Value = 2968.2104
$ws.Range("I100").Value = 2228.2856
$ws.Range("J100").Value = 3399.8333
$ws.Range("K100").Value = 2228.2856
$ws.Range("L100").Value = 3399.8333
$ws.Range("M100").Value = -1687.2856
$ws.Range("N100").Value = -4481.8333
# Row 116
$ws.Range("H116").Value = 8149.875
$ws.Range("I116").Value = 9679.4
$ws.Range("J116").Value = 7454.636
$ws.Range("K116").Value = 9679.4
$ws.Range("L116").Value = 7454.636
$ws.Range("M116").Value = -6237.4
$ws.Range("N116").Value = -14338.636
# Row 132
$ws.Range("H132").Value = 1548.0555
$ws.Range("I132").Value = 1548.0555
$ws.Range("K132").Value = 4644.166499999999
$ws.Range("M132").Value = -2114.166499999999
# Row 133
$ws.Range("H133").Value = 65999.25
$ws.Range("J133").Value = 65999.25
$ws.Range("L133").Value = 65999.25
$ws.Range("N133").Value = -76119.25
# Row 137
$ws.Range("H137").Value = 3126.3928
$ws.Range("I137").Value = 2350.6667
$ws.Range("K137").Value = 7052.000100000001
$ws.Range("M137").Value = -4502.000100000001
# Row 138
$ws.Range("H138").Value = 3966.9678
$ws.Range("J138").Value = 3754.4783
$ws.Range("L138").Value = 11263.4349
$ws.Range("N138").Value = -21543.4349

# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
# Row 8
$ws.Range("H8").Value = 5005002.5
$ws.Range("I8").Value = 5005002.5
$ws.Range("K8").Value = 5005002.5
$ws.Range("M8").Value = -5004858.5
# Row 34
$ws.Range("H34").Value = 209256.25
$ws.Range("I34").Value = 112341.664
$ws.Range("K34").Value = 112341.664
$ws.Range("M34").Value = -112070.664
# Row 61
$ws.Range("H61").Value = 3814.3225
$ws.Range("I61").Value = 2745.92
$ws.Range("K61").Value = 2745.92
$ws.Range("M61").Value = -2533.92
# Row 74
$ws.Range("H74").Value = 23811656
$ws.Range("I74").Value = 37039296
$ws.Range("K74").Value = 37039296
$ws.Range("M74").Value = -37038422
# Row 77
$ws.Range("H77").Value = 23811656
$ws.Range("I77").Value = 37039296
$ws.Range("K77").Value = 185196480
$ws.Range("M77").Value = -185192112
# Row 96
$ws.Range("H96").Value = 16850
$ws.Range("J96").Value = 16850
$ws.Range("L96").Value = 16850
$ws.Range("N96").Value = -22342
# Row 122
$ws.Range("H122").Value = 1747.0416
$ws.Range("I122").Value = 1388.2174
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 4164.6522
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -1714.6522
$ws.Range("N122").Value = -34900
# Row 136
$ws.Range("H136").Value = 3814.3225
$ws.Range("I136").Value = 2745.92
$ws.Range("K136").Value = 8237.76
$ws.Range("M136").Value = -5687.76

# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
# Row 11
$ws.Range("H11").Value = 2280.4
$ws.Range("I11").Value = 1972
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 1972
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -1832
$ws.Range("N11").Value = -3280
# Row 99
$ws.Range("H99").Value = 1692.5
$ws.Range("I99").Value = 1615.8
$ws.Range("J99").Value = 1820.3334
$ws.Range("K99").Value = 1615.8
$ws.Range("L99").Value = 1820.3334
$ws.Range("M99").Value = -117.8
$ws.Range("N99").Value = -4816.3334

# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 25064.857
$ws.Range("I31").Value = 2433.2285
$ws.Range("J31").Value = 81643.92999999999
$ws.Range("K31").Value = 2433.2285
$ws.Range("L31").Value = 81643.92999999999
$ws.Range("M31").Value = -2138.2285
$ws.Range("N31").Value = -82233.92999999999
# Row 34
$ws.Range("H34").Value = 25064.857
$ws.Range("I34").Value = 2433.2285
$ws.Range("J34").Value = 81643.92999999999
$ws.Range("K34").Value = 2433.2285
$ws.Range("L34").Value = 81643.92999999999
$ws.Range("M34").Value = -2231.2285
$ws.Range("N34").Value = -82047.92999999999
# Row 99
$ws.Range("H99").Value = 3057
$ws.Range("J99").Value = 3400
$ws.Range("L99").Value = 3400
$ws.Range("N99").Value = -6396
# Row 106
$ws.Range("H106").Value = 44665.668
$ws.Range("J106").Value = 44665.668
$ws.Range("L106").Value = 44665.668
$ws.Range("N106").Value = -47189.668
# Row 107
$ws.Range("H107").Value = 1499.9412
$ws.Range("I107").Value = 1364.9
$ws.Range("J107").Value = 1692.8572
$ws.Range("K107").Value = 1364.9
$ws.Range("L107").Value = 1692.8572
$ws.Range("M107").Value = 555.0999999999999
$ws.Range("N107").Value = -5532.8572
# Row 126
$ws.Range("H126").Value = 3057
$ws.Range("J126").Value = 3400
$ws.Range("L126").Value = 10200
$ws.Range("N126").Value = -15140
# Row 134
$ws.Range("H134").Value = 2953.8572
$ws.Range("I134").Value = 2412.1538
$ws.Range("J134").Value = 4518.778
$ws.Range("K134").Value = 7236.4614
$ws.Range("L134").Value = 13556.334
$ws.Range("M134").Value = -4701.4614
$ws.Range("N134").Value = -18626.334

# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Range("H5").Value = 8686.8125
$ws.Range("I5").Value = 810.7
$ws.Range("K5").Value = 2432.1
$ws.Range("M5").Value = -2320.1
# Row 105
$ws.Range("H105").Value = 20026
$ws.Range("I105").Value = 20026
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 60078
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -57457
$ws.Range("N105").ClearContents()
# Row 132
$ws.Range("H132").Value = 5275
$ws.Range("J132").Value = 5462.5
$ws.Range("L132").Value = 49162.5
$ws.Range("N132").Value = -54222.5
# Row 135
$ws.Range("H135").Value = 8686.8125
$ws.Range("I135").Value = 810.7
$ws.Range("K135").Value = 7296.3
$ws.Range("M135").Value = -4761.3

# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
# Row 44
$ws.Range("H44").Value = 16119.667
$ws.Range("I44").Value = 13028
$ws.Range("J44").Value = 31578
$ws.Range("K44").Value = 13028
$ws.Range("L44").Value = 31578
$ws.Range("M44").Value = -12432
$ws.Range("N44").Value = -32770
# Row 97
$ws.Range("H97").Value = 2598
$ws.Range("I97").Value = 2295
$ws.Range("J97").Value = 2800
$ws.Range("K97").Value = 2295
$ws.Range("L97").Value = 2800
$ws.Range("M97").Value = -1799
$ws.Range("N97").Value = -3792
# Row 132
$ws.Range("H132").Value = 96668.73
$ws.Range("I132").Value = 168090.5
$ws.Range("J132").Value = 10962.6
$ws.Range("K132").Value = 504271.5
$ws.Range("L132").Value = 32887.8
$ws.Range("M132").Value = -501741.5
$ws.Range("N132").Value = -37947.8

# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Range("H7").Value = 5541.343
$ws.Range("I7").Value = 2752.3462
$ws.Range("K7").Value = 2752.3462
$ws.Range("M7").Value = -2640.3462
# Row 40
$ws.Range("H40").Value = 4510.8
$ws.Range("J40").Value = 14505
$ws.Range("L40").Value = 14505
$ws.Range("N40").Value = -14777
# Row 43
$ws.Range("H43").Value = 28999.5
$ws.Range("I43").Value = 29999
$ws.Range("K43").Value = 29999
$ws.Range("M43").Value = -29806
# Row 55
$ws.Range("H55").Value = 1446.6511
$ws.Range("I55").Value = 1123.12
$ws.Range("J55").Value = 1896
$ws.Range("K55").Value = 1123.12
$ws.Range("L55").Value = 1896
$ws.Range("M55").Value = -950.1199999999999
$ws.Range("N55").Value = -2242
# Row 93
$ws.Range("H93").Value = 3712.875
$ws.Range("I93").Value = 3000.75
$ws.Range("K93").Value = 3000.75
$ws.Range("M93").Value = -1752.75
# Row 126
$ws.Range("H126").Value = 5541.343
$ws.Range("I126").Value = 2752.3462
$ws.Range("K126").Value = 8257.0386
$ws.Range("M126").Value = -5787.0386
# Row 132
$ws.Range("H132").Value = 5224
$ws.Range("I132").Value = 3755.5
$ws.Range("J132").Value = 6398.8
$ws.Range("K132").Value = 11266.5
$ws.Range("L132").Value = 19196.4
$ws.Range("M132").Value = -8736.5
$ws.Range("N132").Value = -24256.4

# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
# Row 29
$ws.Range("H29").Value = 1591918.4
$ws.Range("I29").Value = 787877.5
$ws.Range("K29").Value = 787877.5
$ws.Range("M29").Value = -787587.5
# Row 52
$ws.Range("H52").Value = 19000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15452
# Row 96
$ws.Range("H96").Value = 5724
$ws.Range("J96").Value = 9999
$ws.Range("L96").Value = 9999
$ws.Range("N96").Value = -12745
# Row 107
$ws.Range("H107").Value = 620.125
$ws.Range("J107").Value = 450.66666
$ws.Range("L107").Value = 1351.99998
$ws.Range("N107").Value = -5191.999980000001
# Row 132
$ws.Range("H132").Value = 5322.923
$ws.Range("I132").Value = 4881.5454
$ws.Range("J132").Value = 7750.5
$ws.Range("K132").Value = 14644.6362
$ws.Range("L132").Value = 23251.5
$ws.Range("M132").Value = -12114.6362
$ws.Range("N132").Value = -28311.5
